$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 25 through 35 (11 rows), shifting everything below up.
$ws.Rows("25:35").Delete()

# Update selection/top-left to match the post-edit view.
$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("A25:XFD35").Select()
